$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ E=3; G=2.536023666666666; H=7.608070999999999; I=0.02556389501525096; J=0.02556389501525096; K=3; M=29.57110033333333; N=88.713301; O=0.5311523066901919; P=0.5311523066901919; Q=74.99301029470789; R=674.9370926523709; S=0.01357832180533644; T=0.01357832180533645 }
    3  = @{ E=3; G=2.536023666666666; H=7.608070999999999; I=0.02556389501525096; J=0.02556389501525096; K=3; M=12.180337; N=36.541011; O=0.2187816490048282; P=0.2187816490048282; Q=30.88962289997567; R=278.006606099781; S=0.005592911106422913; T=0.005592911106422914 }
    4  = @{ E=3; G=2.536023666666666; H=7.608070999999999; I=0.02556389501525096; J=0.02556389501525096; K=3; M=13.92204833333333; N=41.76614499999999; O=0.2500660443049799; P=0.2500660443049799; Q=35.30664406181054; R=317.7597965562949; S=0.006392662103491602; T=0.006392662103491603 }
    5  = @{ E=3; G=75.11538433333334; H=225.346153; I=0.7571860721834327; J=0.7571860721834328; K=3; M=29.57110033333333; N=88.713301; O=0.5311523066901919; P=0.5311523066901919; Q=2221.244566697895; R=19991.20110028105; S=0.4021811288339164; T=0.4021811288339164 }
    6  = @{ E=3; G=75.11538433333334; H=225.346153; I=0.7571860721834327; J=0.7571860721834328; K=3; M=12.180337; N=36.541011; O=0.2187816490048282; P=0.2187816490048282; Q=914.9306950645205; R=8234.376255580684; S=0.1656584174757803; T=0.1656584174757803 }
    7  = @{ E=3; G=75.11538433333334; H=225.346153; I=0.7571860721834327; J=0.7571860721834328; K=3; M=13.92204833333333; N=41.76614499999999; O=0.2500660443049799; P=0.2500660443049799; Q=1045.760011265576; R=9411.840101390184; S=0.189346525873736; T=0.189346525873736 }
    8  = @{ E=3; G=21.551928; H=64.655784; I=0.2172500328013163; J=0.2172500328013163; K=3; M=29.57110033333333; N=88.713301; O=0.5311523066901919; P=0.5311523066901919; Q=637.314225264776; R=5735.828027382984; S=0.115392856050939; T=0.115392856050939 }
    9  = @{ E=3; G=21.551928; H=64.655784; I=0.2172500328013163; J=0.2172500328013163; K=3; M=12.180337; N=36.541011; O=0.2187816490048282; P=0.2187816490048282; Q=262.5097460397361; R=2362.587714357624; S=0.04753032042262499; T=0.04753032042262499 }
    10 = @{ E=3; G=21.551928; H=64.655784; I=0.2172500328013163; J=0.2172500328013163; K=3; M=13.92204833333333; N=41.76614499999999; O=0.2500660443049799; P=0.2500660443049799; Q=300.0469832925199; R=2700.42284963268; S=0.0543268563277523; T=0.0543268563277523 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    foreach ($col in $vals.Keys) {
        $ws.Range("$col$row").Value = $vals[$col]
    }
}
